# Update the "Förändrad" (changed) date column (C) from 45174 to 45175
# for all data rows (rows 2 through 27) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
